$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "G2"; Value = 29.20950566666667 }
    @{ Cell = "H2"; Value = 87.628517 }
    @{ Cell = "I2"; Value = 0.01829497698069002 }
    @{ Cell = "J2"; Value = 0.01840828041918582 }
    @{ Cell = "M2"; Value = 2.621797333333333 }
    @{ Cell = "N2"; Value = 7.865392 }
    @{ Cell = "O2"; Value = 0.07867217155043885 }
    @{ Cell = "P2"; Value = 0.07906089226781998 }
    @{ Cell = "Q2"; Value = 76.58140406485155 }
    @{ Cell = "R2"; Value = 689.232636583664 }
    @{ Cell = "S2"; Value = 0.001439305567536175 }
    @{ Cell = "T2"; Value = 0.00145537507505707 }
    @{ Cell = "G3"; Value = 29.20950566666667 }
    @{ Cell = "H3"; Value = 87.628517 }
    @{ Cell = "I3"; Value = 0.01829497698069002 }
    @{ Cell = "J3"; Value = 0.01840828041918582 }
    @{ Cell = "O3"; Value = 0.148308476825081 }
    @{ Cell = "P3"; Value = 0.1490412718702539 }
    @{ Cell = "Q3"; Value = 144.3670762628244 }
    @{ Cell = "R3"; Value = 1299.303686365419 }
    @{ Cell = "S3"; Value = 0.002713300169556057 }
    @{ Cell = "T3"; Value = 0.002743593526619745 }
    @{ Cell = "G4"; Value = 29.20950566666667 }
    @{ Cell = "H4"; Value = 87.628517 }
    @{ Cell = "I4"; Value = 0.01829497698069002 }
    @{ Cell = "J4"; Value = 0.01840828041918582 }
    @{ Cell = "M4"; Value = 11.09754033333333 }
    @{ Cell = "N4"; Value = 33.292621 }
    @{ Cell = "O4"; Value = 0.3330034651388949 }
    @{ Cell = "P4"; Value = 0.3346488416844782 }
    @{ Cell = "Q4"; Value = 324.1536672525619 }
    @{ Cell = "R4"; Value = 2917.383005273057 }
    @{ Cell = "S4"; Value = 0.006092290729206093 }
    @{ Cell = "T4"; Value = 0.006160309719683595 }
    @{ Cell = "G5"; Value = 29.20950566666667 }
    @{ Cell = "H5"; Value = 87.628517 }
    @{ Cell = "I5"; Value = 0.01829497698069002 }
    @{ Cell = "J5"; Value = 0.01840828041918582 }
    @{ Cell = "M5"; Value = 0.4915585 }
    @{ Cell = "N5"; Value = 0.983117 }
    @{ Cell = "O5"; Value = 0.01475017696730553 }
    @{ Cell = "P5"; Value = 0.009882038584175128 }
    @{ Cell = "Q5"; Value = 14.35818079124817 }
    @{ Cell = "R5"; Value = 86.149084747489 }
    @{ Cell = "S5"; Value = 0.0002698541480779588 }
    @{ Cell = "T5"; Value = 0.0001819113373707098 }
    @{ Cell = "G6"; Value = 29.20950566666667 }
    @{ Cell = "H6"; Value = 87.628517 }
    @{ Cell = "I6"; Value = 0.01829497698069002 }
    @{ Cell = "J6"; Value = 0.01840828041918582 }
    @{ Cell = "M6"; Value = 14.172235 }
    @{ Cell = "N6"; Value = 42.516705 }
    @{ Cell = "O6"; Value = 0.4252657095182797 }
    @{ Cell = "P6"; Value = 0.4273669555932728 }
    @{ Cell = "Q6"; Value = 413.9639785418317 }
    @{ Cell = "R6"; Value = 3725.675806876485 }
    @{ Cell = "S6"; Value = 0.007780226366313736 }
    @{ Cell = "T6"; Value = 0.007867090760454699 }
    @{ Cell = "I7"; Value = 0.913374480506715 }
    @{ Cell = "J7"; Value = 0.9190311407684336 }
    @{ Cell = "M7"; Value = 2.621797333333333 }
    @{ Cell = "N7"; Value = 7.865392 }
    @{ Cell = "O7"; Value = 0.07867217155043885 }
    @{ Cell = "P7"; Value = 0.07906089226781998 }
    @{ Cell = "Q7"; Value = 3823.317199471571 }
    @{ Cell = "R7"; Value = 34409.85479524414 }
    @{ Cell = "S7"; Value = 0.07185715382021725 }
    @{ Cell = "T7"; Value = 0.07265942201106483 }
    @{ Cell = "I8"; Value = 0.913374480506715 }
    @{ Cell = "J8"; Value = 0.9190311407684336 }
    @{ Cell = "O8"; Value = 0.148308476825081 }
    @{ Cell = "P8"; Value = 0.1490412718702539 }
    @{ Cell = "S8"; Value = 0.1354611779748505 }
    @{ Cell = "T8"; Value = 0.1369735701084977 }
    @{ Cell = "I9"; Value = 0.913374480506715 }
    @{ Cell = "J9"; Value = 0.9190311407684336 }
    @{ Cell = "M9"; Value = 11.09754033333333 }
    @{ Cell = "N9"; Value = 33.292621 }
    @{ Cell = "O9"; Value = 0.3330034651388949 }
    @{ Cell = "P9"; Value = 0.3346488416844782 }
    @{ Cell = "Q9"; Value = 16183.33205576892 }
    @{ Cell = "R9"; Value = 145649.9885019203 }
    @{ Cell = "S9"; Value = 0.3041568669781741 }
    @{ Cell = "T9"; Value = 0.3075527067301209 }
    @{ Cell = "I10"; Value = 0.913374480506715 }
    @{ Cell = "J10"; Value = 0.9190311407684336 }
    @{ Cell = "M10"; Value = 0.4915585 }
    @{ Cell = "N10"; Value = 0.983117 }
    @{ Cell = "O10"; Value = 0.01475017696730553 }
    @{ Cell = "P10"; Value = 0.009882038584175128 }
    @{ Cell = "Q10"; Value = 716.8304138928281 }
    @{ Cell = "R10"; Value = 4300.982483356969 }
    @{ Cell = "S10"; Value = 0.0134724352248948 }
    @{ Cell = "T10"; Value = 0.009081901193132144 }
    @{ Cell = "I11"; Value = 0.913374480506715 }
    @{ Cell = "J11"; Value = 0.9190311407684336 }
    @{ Cell = "M11"; Value = 14.172235 }
    @{ Cell = "N11"; Value = 42.516705 }
    @{ Cell = "O11"; Value = 0.4252657095182797 }
    @{ Cell = "P11"; Value = 0.4273669555932728 }
    @{ Cell = "Q11"; Value = 20667.10082489963 }
    @{ Cell = "R11"; Value = 186003.9074240967 }
    @{ Cell = "S11"; Value = 0.3884268465085783 }
    @{ Cell = "T11"; Value = 0.392763540725618 }
    @{ Cell = "G12"; Value = 57.98602933333333 }
    @{ Cell = "H12"; Value = 173.958088 }
    @{ Cell = "I12"; Value = 0.03631876156896331 }
    @{ Cell = "J12"; Value = 0.03654368891224535 }
    @{ Cell = "M12"; Value = 2.621797333333333 }
    @{ Cell = "N12"; Value = 7.865392 }
    @{ Cell = "O12"; Value = 0.07867217155043885 }
    @{ Cell = "P12"; Value = 0.07906089226781998 }
    @{ Cell = "Q12"; Value = 152.0276170767218 }
    @{ Cell = "R12"; Value = 1368.248553690496 }
    @{ Cell = "S12"; Value = 0.002857275840652967 }
    @{ Cell = "T12"; Value = 0.002889176652159757 }
    @{ Cell = "G13"; Value = 57.98602933333333 }
    @{ Cell = "H13"; Value = 173.958088 }
    @{ Cell = "I13"; Value = 0.03631876156896331 }
    @{ Cell = "J13"; Value = 0.03654368891224535 }
    @{ Cell = "O13"; Value = 0.148308476825081 }
    @{ Cell = "P13"; Value = 0.1490412718702539 }
    @{ Cell = "Q13"; Value = 286.5941524130906 }
    @{ Cell = "R13"; Value = 2579.347371717816 }
    @{ Cell = "S13"; Value = 0.005386380208466239 }
    @{ Cell = "T13"; Value = 0.005446517874311942 }
    @{ Cell = "G14"; Value = 57.98602933333333 }
    @{ Cell = "H14"; Value = 173.958088 }
    @{ Cell = "I14"; Value = 0.03631876156896331 }
    @{ Cell = "J14"; Value = 0.03654368891224535 }
    @{ Cell = "M14"; Value = 11.09754033333333 }
    @{ Cell = "N14"; Value = 33.292621 }
    @{ Cell = "O14"; Value = 0.3330034651388949 }
    @{ Cell = "P14"; Value = 0.3346488416844782 }
    @{ Cell = "Q14"; Value = 643.5022992965163 }
    @{ Cell = "R14"; Value = 5791.520693668646 }
    @{ Cell = "S14"; Value = 0.01209427345201811 }
    @{ Cell = "T14"; Value = 0.01222930316536081 }
    @{ Cell = "G15"; Value = 57.98602933333333 }
    @{ Cell = "H15"; Value = 173.958088 }
    @{ Cell = "I15"; Value = 0.03631876156896331 }
    @{ Cell = "J15"; Value = 0.03654368891224535 }
    @{ Cell = "M15"; Value = 0.4915585 }
    @{ Cell = "N15"; Value = 0.983117 }
    @{ Cell = "O15"; Value = 0.01475017696730553 }
    @{ Cell = "P15"; Value = 0.009882038584175128 }
    @{ Cell = "Q15"; Value = 28.50352560004933 }
    @{ Cell = "R15"; Value = 171.021153600296 }
    @{ Cell = "S15"; Value = 0.0005357081603755839 }
    @{ Cell = "T15"; Value = 0.0003611261438389013 }
    @{ Cell = "G16"; Value = 57.98602933333333 }
    @{ Cell = "H16"; Value = 173.958088 }
    @{ Cell = "I16"; Value = 0.03631876156896331 }
    @{ Cell = "J16"; Value = 0.03654368891224535 }
    @{ Cell = "M16"; Value = 14.172235 }
    @{ Cell = "N16"; Value = 42.516705 }
    @{ Cell = "O16"; Value = 0.4252657095182797 }
    @{ Cell = "P16"; Value = 0.4273669555932728 }
    @{ Cell = "Q16"; Value = 821.7916344288933 }
    @{ Cell = "R16"; Value = 7396.124709860039 }
    @{ Cell = "S16"; Value = 0.01544512390745041 }
    @{ Cell = "T16"; Value = 0.01561756507657393 }
    @{ Cell = "G17"; Value = 29.481085 }
    @{ Cell = "H17"; Value = 58.96217 }
    @{ Cell = "I17"; Value = 0.01846507700595112 }
    @{ Cell = "J17"; Value = 0.01238628926567028 }
    @{ Cell = "M17"; Value = 2.621797333333333 }
    @{ Cell = "N17"; Value = 7.865392 }
    @{ Cell = "O17"; Value = 0.07867217155043885 }
    @{ Cell = "P17"; Value = 0.07906089226781998 }
    @{ Cell = "Q17"; Value = 77.29343003677333 }
    @{ Cell = "R17"; Value = 463.76058022064 }
    @{ Cell = "S17"; Value = 0.001452687705904251 }
    @{ Cell = "T17"; Value = 0.0009792710812312128 }
    @{ Cell = "G18"; Value = 29.481085 }
    @{ Cell = "H18"; Value = 58.96217 }
    @{ Cell = "I18"; Value = 0.01846507700595112 }
    @{ Cell = "J18"; Value = 0.01238628926567028 }
    @{ Cell = "O18"; Value = 0.148308476825081 }
    @{ Cell = "P18"; Value = 0.1490412718702539 }
    @{ Cell = "Q18"; Value = 145.709348698865 }
    @{ Cell = "R18"; Value = 874.2560921931901 }
    @{ Cell = "S18"; Value = 0.002738527445210439 }
    @{ Cell = "T18"; Value = 0.001846068305908371 }
    @{ Cell = "G19"; Value = 29.481085 }
    @{ Cell = "H19"; Value = 58.96217 }
    @{ Cell = "I19"; Value = 0.01846507700595112 }
    @{ Cell = "J19"; Value = 0.01238628926567028 }
    @{ Cell = "M19"; Value = 11.09754033333333 }
    @{ Cell = "N19"; Value = 33.292621 }
    @{ Cell = "O19"; Value = 0.3330034651388949 }
    @{ Cell = "P19"; Value = 0.3346488416844782 }
    @{ Cell = "Q19"; Value = 327.1675298579283 }
    @{ Cell = "R19"; Value = 1963.00517914757 }
    @{ Cell = "S19"; Value = 0.006148934627038256 }
    @{ Cell = "T19"; Value = 0.004145057355525444 }
    @{ Cell = "G20"; Value = 29.481085 }
    @{ Cell = "H20"; Value = 58.96217 }
    @{ Cell = "I20"; Value = 0.01846507700595112 }
    @{ Cell = "J20"; Value = 0.01238628926567028 }
    @{ Cell = "M20"; Value = 0.4915585 }
    @{ Cell = "N20"; Value = 0.983117 }
    @{ Cell = "O20"; Value = 0.01475017696730553 }
    @{ Cell = "P20"; Value = 0.009882038584175128 }
    @{ Cell = "Q20"; Value = 14.4916779209725 }
    @{ Cell = "R20"; Value = 57.96671168389 }
    @{ Cell = "S20"; Value = 0.0002723631535527032 }
    @{ Cell = "T20"; Value = 0.0001224017884381079 }
    @{ Cell = "G21"; Value = 29.481085 }
    @{ Cell = "H21"; Value = 58.96217 }
    @{ Cell = "I21"; Value = 0.01846507700595112 }
    @{ Cell = "J21"; Value = 0.01238628926567028 }
    @{ Cell = "M21"; Value = 14.172235 }
    @{ Cell = "N21"; Value = 42.516705 }
    @{ Cell = "O21"; Value = 0.4252657095182797 }
    @{ Cell = "P21"; Value = 0.4273669555932728 }
    @{ Cell = "Q21"; Value = 417.812864674975 }
    @{ Cell = "R21"; Value = 2506.87718804985 }
    @{ Cell = "S21"; Value = 0.007852564074245477 }
    @{ Cell = "T21"; Value = 0.005293490734567141 }
    @{ Cell = "G22"; Value = 21.628479 }
    @{ Cell = "H22"; Value = 64.885437 }
    @{ Cell = "I22"; Value = 0.01354670393768061 }
    @{ Cell = "J22"; Value = 0.01363060063446486 }
    @{ Cell = "M22"; Value = 2.621797333333333 }
    @{ Cell = "N22"; Value = 7.865392 }
    @{ Cell = "O22"; Value = 0.07867217155043885 }
    @{ Cell = "P22"; Value = 0.07906089226781998 }
    @{ Cell = "Q22"; Value = 56.70548856625599 }
    @{ Cell = "R22"; Value = 510.349397096304 }
    @{ Cell = "S22"; Value = 0.001065748616128214 }
    @{ Cell = "T22"; Value = 0.001077647448307105 }
    @{ Cell = "G23"; Value = 21.628479 }
    @{ Cell = "H23"; Value = 64.885437 }
    @{ Cell = "I23"; Value = 0.01354670393768061 }
    @{ Cell = "J23"; Value = 0.01363060063446486 }
    @{ Cell = "O23"; Value = 0.148308476825081 }
    @{ Cell = "P23"; Value = 0.1490412718702539 }
    @{ Cell = "Q23"; Value = 106.898086974651 }
    @{ Cell = "R23"; Value = 962.082782771859 }
    @{ Cell = "S23"; Value = 0.002009091026997739 }
    @{ Cell = "T23"; Value = 0.002031522054916133 }
    @{ Cell = "G24"; Value = 21.628479 }
    @{ Cell = "H24"; Value = 64.885437 }
    @{ Cell = "I24"; Value = 0.01354670393768061 }
    @{ Cell = "J24"; Value = 0.01363060063446486 }
    @{ Cell = "M24"; Value = 11.09754033333333 }
    @{ Cell = "N24"; Value = 33.292621 }
    @{ Cell = "O24"; Value = 0.3330034651388949 }
    @{ Cell = "P24"; Value = 0.3346488416844782 }
    @{ Cell = "Q24"; Value = 240.022918051153 }
    @{ Cell = "R24"; Value = 2160.206262460377 }
    @{ Cell = "S24"; Value = 0.004511099352458355 }
    @{ Cell = "T24"; Value = 0.004561464713787379 }
    @{ Cell = "G25"; Value = 21.628479 }
    @{ Cell = "H25"; Value = 64.885437 }
    @{ Cell = "I25"; Value = 0.01354670393768061 }
    @{ Cell = "J25"; Value = 0.01363060063446486 }
    @{ Cell = "M25"; Value = 0.4915585 }
    @{ Cell = "N25"; Value = 0.983117 }
    @{ Cell = "O25"; Value = 0.01475017696730553 }
    @{ Cell = "P25"; Value = 0.009882038584175128 }
    @{ Cell = "Q25"; Value = 10.6316626945215 }
    @{ Cell = "R25"; Value = 63.789976167129 }
    @{ Cell = "S25"; Value = 0.0001998162804044837 }
    @{ Cell = "T25"; Value = 0.0001346981213952637 }
    @{ Cell = "G26"; Value = 21.628479 }
    @{ Cell = "H26"; Value = 64.885437 }
    @{ Cell = "I26"; Value = 0.01354670393768061 }
    @{ Cell = "J26"; Value = 0.01363060063446486 }
    @{ Cell = "M26"; Value = 14.172235 }
    @{ Cell = "N26"; Value = 42.516705 }
    @{ Cell = "O26"; Value = 0.4252657095182797 }
    @{ Cell = "P26"; Value = 0.4273669555932728 }
    @{ Cell = "Q26"; Value = 306.523887080565 }
    @{ Cell = "R26"; Value = 2758.714983725085 }
    @{ Cell = "S26"; Value = 0.005760948661691818 }
    @{ Cell = "T26"; Value = 0.005825268296058981 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

